# Update cryptocurrency price (column D) and 1h volume change (column E) values
# for rows 2-51 on the active worksheet, per the latest data refresh.
# Column D cells are stored as plain text (not numbers), so force a text
# number format before assigning to stop Excel from auto-converting values
# such as "1.0000" or "0.9994" into numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.506.86"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.727.56"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.91"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4805"
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06218"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.726.06"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07154"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6178"
$ws.Range("E13").Value = "  +4.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.521"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.15"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.515.79"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.0000"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006929"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.66"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.948.56"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.530"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.955"
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.275"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.16"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.32"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.407"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.65"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.979"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08027"
$ws.Range("E31").Value = "  +3.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.707"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9993"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.614"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6357"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9907"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9326"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.097"
$ws.Range("E39").Value = "  +9.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.416"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.63"
$ws.Range("E41").Value = "  -9.79%  "
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01503"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.650"
$ws.Range("E44").Value = "  +6.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3914"
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.911"
$ws.Range("E46").Value = "  +10.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1185"
$ws.Range("E47").Value = "  +2.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05332"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.97"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.862"
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.265"
$ws.Range("E51").Value = "  +3.55%  "

Write-Host "Updated cryptos list"
